$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G5").ClearContents()
